$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 (I0) and J1 (IF), copying the header style from H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# Populate I0 / IF data for rows 2-77
$i0 = @{
    2 = 8
    3 = 9
    4 = 8
    5 = 8
    6 = 9
    7 = 8
    8 = 9
    9 = 9
    10 = 9
    11 = 9
    12 = 8
    13 = 9
    14 = 9
    15 = 10
    16 = 9
    17 = 9
    18 = 9
    19 = 8
    20 = 9
    21 = 9
    22 = 9
    23 = 9
    24 = 8
    25 = 9
    26 = 9
    27 = 9
    28 = 9
    29 = 9
    30 = 9
    31 = 8
    32 = 9
    33 = 10
    34 = 9
    35 = 9
    36 = 9
    37 = 9
    38 = 9
    39 = 8
    40 = 9
    41 = 9
    42 = 9
    43 = 9
    44 = 9
    45 = 9
    46 = 9
    47 = 9
    48 = 9
    49 = 9
    50 = 9
    51 = 9
    52 = 9
    53 = 9
    54 = 8
    55 = 9
    56 = 9
    57 = 8
    58 = 9
    59 = 8
    60 = 9
    61 = 10
    62 = 9
    63 = 9
    64 = 9
    65 = 9
    66 = 8
    67 = 8
    68 = 9
    69 = 9
    70 = 8
    71 = 9
    72 = 9
    73 = 5
    74 = 9
    75 = 8
    76 = 4
    77 = 3
}
$if = @{
    2 = 8
    3 = 9
    4 = 8
    5 = 8
    6 = 9
    7 = 9
    8 = 9
    9 = 9
    10 = 9
    11 = 9
    12 = 8
    13 = 9
    14 = 9
    15 = 10
    16 = 9
    17 = 9
    18 = 9
    19 = 9
    20 = 9
    21 = 9
    22 = 10
    23 = 9
    24 = 9
    25 = 9
    26 = 9
    27 = 9
    28 = 9
    29 = 9
    30 = 9
    31 = 8
    32 = 9
    33 = 10
    34 = 9
    35 = 9
    36 = 9
    37 = 9
    38 = 9
    39 = 8
    40 = 9
    41 = 9
    42 = 9
    43 = 9
    44 = 9
    45 = 9
    46 = 9
    47 = 9
    48 = 9
    49 = 9
    50 = 9
    51 = 9
    52 = 9
    53 = 9
    54 = 9
    55 = 9
    56 = 9
    57 = 9
    58 = 9
    59 = 9
    60 = 9
    61 = 10
    62 = 9
    63 = 9
    64 = 9
    65 = 9
    66 = 8
    67 = 8
    68 = 9
    69 = 9
    70 = 9
    71 = 9
    72 = 9
    73 = 5
    74 = 9
    75 = 8
    76 = 4
    77 = 3
}

foreach ($r in 2..77) {
    $ws.Cells.Item($r, 9).Value = $i0[$r]
    $ws.Cells.Item($r, 10).Value = $if[$r]
}